# Updated the screen recording data in the "Document to Download" log sheet.
# Rows 2-6 (columns E:G = Start time / End time / Time taken) get refreshed
# timestamps from the new recording run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2023-07-12 01:24:03"
$ws.Range("F2").Value = "2023-07-12 01:25:36"
$ws.Range("G2").Value = "00:01:33"

$ws.Range("E3").Value = "2023-07-12 01:25:39"
$ws.Range("F3").Value = "2023-07-12 01:27:12"
$ws.Range("G3").Value = "00:01:33"

$ws.Range("E4").Value = "2023-07-12 01:27:14"
$ws.Range("F4").Value = "2023-07-12 01:28:46"
$ws.Range("G4").Value = "00:01:32"

$ws.Range("E5").Value = "2023-07-12 01:28:49"
$ws.Range("F5").Value = "2023-07-12 01:30:20"
$ws.Range("G5").Value = "00:01:31"

$ws.Range("E6").Value = "2023-07-12 01:30:23"
$ws.Range("F6").Value = "2023-07-12 01:39:36"
$ws.Range("G6").Value = "00:09:13"
